$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "20.00"
# or "315.94" are not auto-coerced into numbers by COM assignment,
# matching the source data which stores these as strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.478.86'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '1.863.36'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').Value = '315.94'
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = '1.003'
$ws.Range('D7').Value = '0.4674'
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').Value = '0.3724'
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('D9').Value = '0.07377'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').Value = '0.8883'
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('D12').Value = '20.00'
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('D13').Value = '1.844.39'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = '5.425'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').Value = '6.593'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '92.68'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '0.000008930'
$ws.Range('E18').Value = '  +3.64%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').Value = '14.92'
$ws.Range('E20').Value = '  +3.31%  '
$ws.Range('D21').Value = '27.514.84'
$ws.Range('E21').Value = '  +2.50%  '
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '10.56'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '2.091.97'
$ws.Range('E24').Value = '  +7.64%  '
$ws.Range('D25').Value = '153.24'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('D26').Value = '1.880'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('E27').Value = '  +2.21%  '
$ws.Range('D28').Value = '2.089'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').Value = '117.05'
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('D31').Value = '0.08911'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = '3.027'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7556'
$ws.Range('E33').Value = '  +5.70%  '
$ws.Range('D34').Value = '1.164'
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('D35').Value = '4.488'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').Value = '2.640'
$ws.Range('E36').Value = '  +9.81%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01970'
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.083'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('D39').Value = '0.05283'
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('D40').Value = '2.996'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').Value = '7.180'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = '0.5213'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('D43').Value = '0.1645'
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').Value = '8.357'
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('D45').Value = '0.4877'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('D46').Value = '10.33'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('D47').Value = '1.003'
$ws.Range('D48').Value = '103.84'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').Value = '1.659'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').Value = '0.06262'
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '66.03'
$ws.Range('E51').Value = '  +3.07%  '

# Restore the column's number format back to General to match the
# original workbook formatting (the underlying cells are still text).
$ws.Range("D2:D51").NumberFormat = "General"
